$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Swap Noruega / Suecia order (and update Noruega's stats) ---
$ws.Range("A23").Value = "Noruega"
$ws.Range("B23").Value = 5071
$ws.Range("C23").Value = 194
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = 5012
$ws.Range("F23").Value = 105
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 46

$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 4947
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 103
$ws.Range("E24").Value = 4605
$ws.Range("F24").Value = 393
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 239

# --- Update last-updated timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 13:50"

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 78983
$ws.Range("C8").Value = 1002
$ws.Range("E8").Value = 58860
$ws.Range("G8").Value = 17
$ws.Range("H8").Value = 948

# --- Suiza (row 12) ---
$ws.Range("B12").Value = 18267
$ws.Range("C12").Value = 499
$ws.Range("E12").Value = 13749

# --- Brasil (row 20) ---
$ws.Range("E20").Value = 6559
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 245

# --- Israel (row 21) ---
$ws.Range("E21").Value = 5889
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 33

# --- Eslovenia (row 55) ---
$ws.Range("D55").Value = 70
$ws.Range("E55").Value = 810
$ws.Range("F55").Value = 31
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 17

# --- Sri Lanka (row 108) ---
$ws.Range("B108").Value = 150
$ws.Range("C108").Value = 4
$ws.Range("E108").Value = 126
